$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 8-12: match rows were reordered (source data refresh); update the
# B/F/G/H/I/J..AB fields (E = HomeTeam stays put per row) to their new values ---
# Row 8
$ws.Range("B8").Value2 = 6627737
$ws.Range("F8").Value2 = "Lokomotiv 1929 Sofia"
$ws.Range("G8").Value2 = 2
$ws.Range("H8").Value2 = 0
$ws.Range("I8").Value2 = "H"
$ws.Range("J8").Value2 = 1.5
$ws.Range("K8").Value2 = 3.75
$ws.Range("L8").Value2 = 6.5
$ws.Range("M8").Value2 = 1.444
$ws.Range("N8").Value2 = 4.333
$ws.Range("O8").Value2 = 8
$ws.Range("P8").Value2 = -1.25
$ws.Range("Q8").Value2 = 2
$ws.Range("R8").Value2 = 1.85
$ws.Range("S8").Value2 = 2.25
$ws.Range("T8").Value2 = 1.875
$ws.Range("U8").Value2 = 1.975
$ws.Range("V8").Value2 = 0.444
$ws.Range("X8").Value2 = -1
$ws.Range("Y8").Value2 = 1
$ws.Range("Z8").Value2 = -1
$ws.Range("AB8").Value2 = 0.4875

# Row 9
$ws.Range("B9").Value2 = 6627736
$ws.Range("F9").Value2 = "Arda Kardzhali"
$ws.Range("G9").Value2 = 0
$ws.Range("H9").Value2 = 3
$ws.Range("I9").Value2 = "A"
$ws.Range("J9").Value2 = 5.25
$ws.Range("K9").Value2 = 3.6
$ws.Range("L9").Value2 = 1.571
$ws.Range("M9").Value2 = 26
$ws.Range("N9").Value2 = 11
$ws.Range("O9").Value2 = 1.083
$ws.Range("P9").Value2 = 2.5
$ws.Range("Q9").Value2 = 1.825
$ws.Range("R9").Value2 = 2.025
$ws.Range("S9").Value2 = 3.25
$ws.Range("T9").Value2 = 2
$ws.Range("U9").Value2 = 1.85
$ws.Range("V9").Value2 = -1
$ws.Range("X9").Value2 = 0.08299999999999996
$ws.Range("Y9").Value2 = -1
$ws.Range("Z9").Value2 = 1.025
$ws.Range("AB9").Value2 = 0.425

# Row 10
$ws.Range("B10").Value2 = 6627725
$ws.Range("F10").Value2 = "CSKA Sofia"
$ws.Range("H10").Value2 = 2
$ws.Range("I10").Value2 = "A"
$ws.Range("J10").Value2 = 2.625
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 2.6
$ws.Range("M10").Value2 = 2.55
$ws.Range("N10").Value2 = 3.3
$ws.Range("O10").Value2 = 2.8
$ws.Range("P10").Value2 = 0
$ws.Range("Q10").Value2 = 1.825
$ws.Range("R10").Value2 = 2.025
$ws.Range("S10").Value2 = 2.25
$ws.Range("W10").Value2 = -1
$ws.Range("X10").Value2 = 1.8
$ws.Range("Z10").Value2 = 1.025
$ws.Range("AA10").Value2 = -0.5
$ws.Range("AB10").Value2 = 0.5

# Row 11
$ws.Range("B11").Value2 = 6627724
$ws.Range("F11").Value2 = "Lokomotiv Plovdiv"
$ws.Range("H11").Value2 = 0
$ws.Range("I11").Value2 = "D"
$ws.Range("J11").Value2 = 1.5
$ws.Range("K11").Value2 = 3.8
$ws.Range("L11").Value2 = 6
$ws.Range("M11").Value2 = 1.45
$ws.Range("N11").Value2 = 4.2
$ws.Range("O11").Value2 = 8
$ws.Range("P11").Value2 = -1.25
$ws.Range("Q11").Value2 = 2.025
$ws.Range("R11").Value2 = 1.825
$ws.Range("S11").Value2 = 2.5
$ws.Range("T11").Value2 = 1.85
$ws.Range("U11").Value2 = 2
$ws.Range("W11").Value2 = 3.2
$ws.Range("X11").Value2 = -1
$ws.Range("Y11").Value2 = -1
$ws.Range("Z11").Value2 = 0.825
$ws.Range("AB11").Value2 = 1

# Row 12
$ws.Range("B12").Value2 = 6627290
$ws.Range("F12").Value2 = "Ludogorets Razgrad"
$ws.Range("H12").Value2 = 1
$ws.Range("J12").Value2 = 6
$ws.Range("K12").Value2 = 4
$ws.Range("L12").Value2 = 1.45
$ws.Range("M12").Value2 = 6
$ws.Range("N12").Value2 = 4.333
$ws.Range("O12").Value2 = 1.55
$ws.Range("P12").Value2 = 1
$ws.Range("Q12").Value2 = 2
$ws.Range("R12").Value2 = 1.85
$ws.Range("S12").Value2 = 2.75
$ws.Range("T12").Value2 = 1.825
$ws.Range("U12").Value2 = 2.025
$ws.Range("X12").Value2 = 0.55
$ws.Range("Y12").Value2 = 0
$ws.Range("Z12").Value2 = 0
$ws.Range("AA12").Value2 = -1
$ws.Range("AB12").Value2 = 1.025

# --- Rows 264-271: updated odds values ---
# Row 264
$ws.Range("M264").Value2 = 2.4
$ws.Range("N264").Value2 = 3
$ws.Range("O264").Value2 = 3.2
$ws.Range("Q264").Value2 = 2.05
$ws.Range("R264").Value2 = 1.8
$ws.Range("T264").Value2 = 2.1
$ws.Range("U264").Value2 = 1.775

# Row 265
$ws.Range("M265").Value2 = 2.55
$ws.Range("O265").Value2 = 2.75

# Row 266
$ws.Range("M266").Value2 = 2.25
$ws.Range("N266").Value2 = 3.25
$ws.Range("O266").Value2 = 3.2
$ws.Range("Q266").Value2 = 1.975
$ws.Range("R266").Value2 = 1.875
$ws.Range("S266").Value2 = 2.5
$ws.Range("T266").Value2 = 2.05
$ws.Range("U266").Value2 = 1.8

# Row 267
$ws.Range("M267").Value2 = 1.65
$ws.Range("O267").Value2 = 5.5

# Row 268
$ws.Range("M268").Value2 = 1.909
$ws.Range("N268").Value2 = 3.3
$ws.Range("O268").Value2 = 4.2
$ws.Range("Q268").Value2 = 1.925
$ws.Range("R268").Value2 = 1.925

# Row 269
$ws.Range("T269").Value2 = 1.95
$ws.Range("U269").Value2 = 1.9

# Row 270
$ws.Range("M270").Value2 = 4.2
$ws.Range("N270").Value2 = 3.3
$ws.Range("O270").Value2 = 1.909
$ws.Range("Q270").Value2 = 1.925
$ws.Range("R270").Value2 = 1.925
$ws.Range("T270").Value2 = 2.025
$ws.Range("U270").Value2 = 1.825

# Row 271
$ws.Range("T271").Value2 = 1.95
$ws.Range("U271").Value2 = 1.9
